# Fan curve calculations for the WindowAC weekly validation case.
# The validation run uses a different maximum flow rate, so update the
# "Max flowrate" input cell (E1) on Sheet1; every other cell in the sheet
# (Max power, FF/PLF/Power/Pressure curve, etc.) derives from E1 via
# formulas and recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("E1").Formula = "=0.43157"

# Restore the view state (zoom + active selection) the workbook was saved with.
[void]$ws.Range("D8:D18").Select()
$excel.ActiveWindow.Zoom = 175
